$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("Z5").Value = 36294
$ws.Range("AA5").Value = 39893

# Row 6
$ws.Range("I6").Value = 1
$ws.Range("N6").Value = 1
$ws.Range("R6").Value = 1

# Row 7
$ws.Range("X7").Value = 1
$ws.Range("AI7").Value = 1
$ws.Range("AK7").Value = 1

# Row 8
$ws.Range("AX8").Value = 1
$ws.Range("AY8").Value = 51

# Row 9
$ws.Range("E9").Value = 2364
$ws.Range("F9").Value = 2145
$ws.Range("K9").Value = 13373
$ws.Range("L9").Value = 11422
$ws.Range("O9").Value = 15052
$ws.Range("P9").Value = 14567
$ws.Range("S9").Value = 11370
$ws.Range("T9").Value = 11573
$ws.Range("U9").Value = 1
$ws.Range("V9").Value = 2775
$ws.Range("W9").Value = 2786
$ws.Range("AN9").Value = 1

# Row 10
$ws.Range("H10").Value = 9
$ws.Range("AN10").Value = 1

# Row 12
$ws.Range("AB12").Value = 2572
$ws.Range("AI12").ClearContents()

# Row 13
$ws.Range("Z13").Value = 35045

# Row 14
$ws.Range("J14").Value = 1
$ws.Range("AO14").Value = 1

# Row 15
$ws.Range("X15").Value = 1

# Row 26
$ws.Range("Z26").Value = 23147
$ws.Range("AM26").Value = 0

# Row 41
$ws.Range("J41").Value = 0
$ws.Range("N41").Value = 0
$ws.Range("R41").Value = 0
$ws.Range("Y41").Value = 0
$ws.Range("AE41").Value = 0
